$wb = $excel.ActiveWorkbook

# --- Sheet1 "Canada": update D7 and append row 8 ---
$ws1 = $wb.Worksheets.Item("Canada")

$ws1.Range("D7").Value = 2452.6

$ws1.Range("A8").Value = 44013
$ws1.Range("A8").NumberFormat = $ws1.Range("A7").NumberFormat
$ws1.Range("B8").Value = "Canada"
$ws1.Range("B8").NumberFormat = $ws1.Range("B7").NumberFormat
$ws1.Range("C8").Value = 89.3
$ws1.Range("D8").Value = 2183.6

# --- Sheet2 "Province": update D52:D61 and append rows 62-71 ---
$ws2 = $wb.Worksheets.Item("Province")

$ws2.Range("D52").Value = 40.2
$ws2.Range("D53").Value = 13.2
$ws2.Range("D54").Value = 64.6
$ws2.Range("D55").Value = 38.7
$ws2.Range("D56").Value = 485.3
$ws2.Range("D57").Value = 943.3
$ws2.Range("D58").Value = 69.1
$ws2.Range("D59").Value = 70.3
$ws2.Range("D60").Value = 382.5
$ws2.Range("D61").Value = 345.5

$provinces = @(
    @{ Row = 62; Name = "Newfoundland & Labrador"; C = 18.4; D = 38.6 },
    @{ Row = 63; Name = "Prince Edward Island";     C = 37.5; D = 9.9 },
    @{ Row = 64; Name = "Nova Scotia";               C = 41.6; D = 52.8 },
    @{ Row = 65; Name = "New Brunswick";             C = 15.2; D = 38 },
    @{ Row = 66; Name = "Quebec";                     C = 94.7; D = 436.8 },
    @{ Row = 67; Name = "Ontario";                    C = 98.2; D = 886.6 },
    @{ Row = 68; Name = "Manitoba";                   C = 41.5; D = 56.3 },
    @{ Row = 69; Name = "Saskatchewan";               C = 59.9; D = 53.1 },
    @{ Row = 70; Name = "Alberta";                    C = 77.2; D = 316.1 },
    @{ Row = 71; Name = "British Columbia";           C = 145.9; D = 295.3 }
)

foreach ($p in $provinces) {
    $r = $p.Row
    $ws2.Range("A$r").Value = 44013
    $ws2.Range("A$r").NumberFormat = $ws2.Range("A61").NumberFormat
    $ws2.Range("B$r").Value = $p.Name
    $ws2.Range("C$r").Value = $p.C
    $ws2.Range("D$r").Value = $p.D
}

# First row of the new date-group (like B42 / B52) also carries the date
# number-format style on column B.
$ws2.Range("B62").NumberFormat = $ws2.Range("B52").NumberFormat

# --- Restore/update cell selections on each sheet ---
$ws1.Range("C9").Select()
$ws2.Activate()
$ws2.Range("C63").Select()
